$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.904.84'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '2.823.98'
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '354.30'
$ws.Range("E5").Value = '  +6.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.44'
$ws.Range("E6").Value = '  -3.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.563'
$ws.Range("E7").Value = '  +4.61%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +4.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.61'
$ws.Range("E10").Value = '  -1.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0853'
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("E13").Value = '  -1.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.77'
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("D15").Value = '3.268.28'
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").Value = '2.826.27'
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.889'
$ws.Range("D18").Value = '51.795.19'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.50'
$ws.Range("E19").Value = '  +9.00%  '
$ws.Range("E20").Value = '  -2.87%  '
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("D22").Value = '0.0₃0989'
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '270.34'
$ws.Range("E23").Value = '  -3.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.79'
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("E25").Value = '  +3.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.81'
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("E30").Value = '  -1.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '50.76'
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.04'
$ws.Range("E32").Value = '  -3.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0456'
$ws.Range("E33").Value = '  +26.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.87'
$ws.Range("E34").Value = '  +5.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.32'
$ws.Range("E35").Value = '  +5.82%  '
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  -1.82%  '
$ws.Range("E39").Value = '  -1.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.25'
$ws.Range("E40").Value = '  -4.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.95'
$ws.Range("E41").Value = '  +2.70%  '
$ws.Range("E42").Value = '  +2.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '126.42'
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.31'
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").Value = '2.095.18'
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.37'
$ws.Range("E47").Value = '  +0.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.28'
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.70'
$ws.Range("E49").Value = '  +2.99%  '
$ws.Range("E50").Value = '  +6.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '61.04'
$ws.Range("E51").Value = '  +0.31%  '
